# "created regular proposal page obj" -- populate the (until now empty)
# Sheet2 with the same EstimatedTrialLength / SafetyConcerns / ... proposal
# fields that already live on Sheet1, then leave the workbook pointed at
# the new sheet/selection the way the author had it when they saved.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet2: header row + one data row (mirrors Sheet1's proposal fields,
#     minus the TreatmentArm1/ActiveCare/FollowUp trial columns) -------
$ws2.Range("A1").Value = "EstimatedTrialLength"
$ws2.Range("B1").Value = "SafetyConcerns"
$ws2.Range("C1").Value = "HealthImpact"
$ws2.Range("D1").Value = "ReductionInHealthCare"
$ws2.Range("E1").Value = "Upload1"
$ws2.Range("F1").Value = "Upload2"
$ws2.Range("G1").Value = "Upload3"
$ws2.Range("H1").Value = "BioSketch"

$ws2.Range("A2").Value = 12
$ws2.Range("B2").Value = "Testing concers"
$ws2.Range("C2").Value = "Testing health"
$ws2.Range("D2").Value = "Testing healthcare"
$ws2.Range("E2").Value = "C:\Users\Biswajit.Ghosh\Desktop\Test_Data\TEST DATA UPLOAD\Sample Files\Files\MSWord\Sample MSWord file -1.docx"
$ws2.Range("F2").Value = "C:\Users\Biswajit.Ghosh\Desktop\Test_Data\TEST DATA UPLOAD\Sample Files\Files\MSWord\Sample MSWord file -2.docx"
$ws2.Range("G2").Value = "C:\Users\Biswajit.Ghosh\Desktop\Test_Data\TEST DATA UPLOAD\Sample Files\Files\MSWord\Sample MSWord file -3.docx"
$ws2.Range("H2").Value = "C:\Users\Biswajit.Ghosh\Desktop\Test_Data\TEST DATA UPLOAD\Sample Files\Files\MSWord\Sample MSWord file -5.docx"

# Match the text-style formatting (s="1") Sheet1 uses for its header/data
# cells. Applied after the values are written so the numeric A2 stays a
# real number instead of being coerced to text.
$ws2.Range("A1:H2").NumberFormat = "@"

# Widen the columns that now hold real content.
$ws2.Columns("A").ColumnWidth = 20
$ws2.Columns("B").ColumnWidth = 14.86

# --- View/selection bookkeeping -----------------------------------
# Sheet1 is no longer the sheet the user left open; its old single-cell
# selection is replaced by the full block that was just mirrored onto
# Sheet2.
$ws1.Range("A1:K2").Select()

# Sheet2 becomes the active/visible tab, selected over the data just
# entered, matching the "regular proposal page" now being the focus.
$ws2.Activate()
$ws2.Range("A1:H2").Select()
